# 4.1.3. Функции продукта — product functions sheet update
#
# 1) Clarify system definitions:
#    - A10: "Система отслеживания доступных парковочных мест"
#           -> "Система выбора доступных парковочных мест"
#    - B17: "Передача данных для формирование отчетов "
#           -> "Передача данных для формирования отчетов " (typo fix)
#
# 2) Remove the obsolete "vehicle database" subsystem row (old row 26):
#    A="Формирование баз данных клиентов и транспортных средств",
#    B="Формирование и ведение базы данных о транспортных средствах в паркинге", ...
#    All subsequent rows shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("A10").Value = "Система выбора доступных парковочных мест"
$ws.Range("B17").Value = "Передача данных для формирования отчетов "

$ws.Rows.Item(26).Delete()

$ws.Rows.Item(10).RowHeight = 58.5

$ws.Range("B27").Select()
